# Updated BEA and FRB pensions data
# Applies the January 26, 2018 data refresh to the BenefitEntitlements sheet:
#  - Revises 2014 (col Q) and 2015 (col R) estimates for every state row (5-56)
#  - Updates the "Last updated" footnote text in A58
#  - Leaves the whole-sheet selection state as left by Excel after the refresh

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BenefitEntitlements")

$newValues = @{
    5  = @(5236681000, 5488753000)
    6  = @(53607519, 55745179)
    7  = @(25178187, 26587060)
    8  = @(74421757, 77812569)
    9  = @(31796144, 33378739)
    10 = @(1008912611, 1056303744)
    11 = @(84991581, 88232226)
    12 = @(81840987, 86732513)
    13 = @(10998218, 11907225)
    14 = @(7105423, 7739007)
    15 = @(216593937, 225416743)
    16 = @(137533742, 143812406)
    17 = @(24121242, 25359278)
    18 = @(16118708, 16985446)
    19 = @(343600558, 362706575)
    20 = @(40511789, 43612543)
    21 = @(35092258, 36862411)
    22 = @(27102924, 28176882)
    23 = @(65223345, 68202895)
    24 = @(76635693, 79300577)
    25 = @(15911762, 16349726)
    26 = @(92776849, 98070702)
    27 = @(145260243, 154742932)
    28 = @(129279963, 130938364)
    29 = @(83184366, 88183601)
    30 = @(43855673, 46698266)
    31 = @(101514442, 105856355)
    32 = @(14184744, 14798444)
    33 = @(19333822, 20230963)
    34 = @(49813276, 52682198)
    35 = @(12425976, 13219536)
    36 = @(174342907, 180536978)
    37 = @(40624655, 43118025)
    38 = @(520823150, 545426907)
    39 = @(100128113, 105538049)
    40 = @(7030690, 7622508)
    41 = @(252814009, 262392833)
    42 = @(41166623, 43067363)
    43 = @(79812026, 88548603)
    44 = @(192770675, 200412357)
    45 = @(16857972, 17603873)
    46 = @(53662125, 55508614)
    47 = @(15814532, 16714953)
    48 = @(54537785, 56653832)
    49 = @(269152554, 288011937)
    50 = @(30341447, 31658293)
    51 = @(6919344, 7307974)
    52 = @(105319536, 109867581)
    53 = @(95394242, 97734817)
    54 = @(18421712, 19207389)
    55 = @(81816026, 84623743)
    56 = @(10003138, 10551266)
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Cells.Item($row, 17).Value = $vals[0]  # column Q = 2014
    $ws.Cells.Item($row, 18).Value = $vals[1]  # column R = 2015
}

# Refresh the "Last updated" footnote (row 58, column A)
$ws.Cells.Item(58, 1).Value = "Last updated: January 26, 2018--revised estimates for 2014 and 2015."

# Select the whole sheet, matching the post-refresh selection state
$ws.Cells.Select() | Out-Null
